# Auto-generated by analysis of the commit diff.
# Updates per-sheet "market price" snapshot columns (H-N) to match the
# scheduled runner's refreshed values. No formulas are involved — every
# target cell stores a plain numeric literal.

$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
$updates = @{
    "H70" = 2591.2693
    "I70" = 1630
    "K70" = 4890
    "M70" = -4620
    "H73" = 2591.2693
    "I73" = 1630
    "K73" = 4890
    "M73" = -3954
    "H112" = 1820.0286
    "J112" = 1870.5333
    "L112" = 5611.5999
    "N112" = -7827.5999
    "H137" = 2320.6365
    "I137" = 2157.8723
    "K137" = 6473.6169
    "M137" = -3923.6169
}
foreach ($cellRef in $updates.Keys) {
    $ws.Range($cellRef).Value = $updates[$cellRef]
}

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
$updates = @{
    "H61" = 2461.2058
    "I61" = 2298.8215
    "K61" = 2298.8215
    "M61" = -2086.8215
    "H97" = 1108.0212
    "I97" = 673.5641000000001
    "J97" = 3226
    "K97" = 673.5641000000001
    "L97" = 3226
    "M97" = -177.5641000000001
    "N97" = -4218
    "H122" = 1855.421
    "I122" = 1288.9231
    "J122" = 3082.8333
    "K122" = 3866.7693
    "L122" = 9248.499899999999
    "M122" = -1416.7693
    "N122" = -14148.4999
    "H136" = 2461.2058
    "I136" = 2298.8215
    "K136" = 6896.4645
    "M136" = -4346.4645
}
foreach ($cellRef in $updates.Keys) {
    $ws.Range($cellRef).Value = $updates[$cellRef]
}

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
$updates = @{
    "H20" = 1352.2727
    "I20" = 1416.0714
    "J20" = 1240.625
    "K20" = 1416.0714
    "L20" = 1240.625
    "M20" = -1169.0714
    "N20" = -1734.625
    "H99" = 28732.666
    "I99" = 38299.09
    "K99" = 38299.09
    "M99" = -36801.09
    "H126" = 0
    "J126" = 0
    "L126" = 0
    "H134" = 1771.9445
    "I134" = 1281.1957
    "K134" = 3843.5871
    "M134" = -1308.5871
}
foreach ($cellRef in $updates.Keys) {
    $ws.Range($cellRef).Value = $updates[$cellRef]
}
$ws.Range("N126").ClearContents()

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
$updates = @{
    "H31" = 2248
    "I31" = 2123.2693
    "J31" = 2653.375
    "K31" = 2123.2693
    "L31" = 2653.375
    "M31" = -1828.2693
    "N31" = -3243.375
    "H34" = 2248
    "I34" = 2123.2693
    "J34" = 2653.375
    "K34" = 2123.2693
    "L34" = 2653.375
    "M34" = -1921.2693
    "N34" = -3057.375
    "H58" = 2426.375
    "I58" = 1455.2354
    "K58" = 1455.2354
    "M58" = -1252.2354
    "H122" = 542185.4
    "I122" = 930683.8
    "K122" = 2792051.4
    "M122" = -2789601.4
    "H132" = 1997
    "I132" = 1876.88
    "K132" = 5630.64
    "M132" = -3100.64
    "H134" = 3274.7144
    "I134" = 2837.7222
    "J134" = 4061.3
    "K134" = 8513.1666
    "L134" = 12183.9
    "M134" = -5978.1666
    "N134" = -17253.9
    "H136" = 2426.375
    "I136" = 1455.2354
    "K136" = 4365.706200000001
    "M136" = -1815.706200000001
}
foreach ($cellRef in $updates.Keys) {
    $ws.Range($cellRef).Value = $updates[$cellRef]
}

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
$updates = @{
    "H36" = 1395.7778
    "I36" = 1395.7778
    "K36" = 4187.3334
    "M36" = -4018.3334
    "H86" = 974.5
    "I86" = 299.33334
    "J86" = 3000
    "K86" = 898.0000200000001
    "L86" = 9000
    "M86" = 287.9999799999999
    "N86" = -11372
    "H89" = 974.5
    "I89" = 299.33334
    "J89" = 3000
    "K89" = 2694.00006
    "L89" = 27000
    "M89" = 3233.99994
    "N89" = -38856
    "H132" = 2099.6667
    "I132" = 2099.6667
    "K132" = 18897.0003
    "M132" = -16367.0003
}
foreach ($cellRef in $updates.Keys) {
    $ws.Range($cellRef).Value = $updates[$cellRef]
}

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
$updates = @{
    "H86" = 119999.25
    "J86" = 119999.25
    "L86" = 119999.25
    "N86" = -122371.25
    "H89" = 119999.25
    "J89" = 119999.25
    "L89" = 359997.75
    "N89" = -371853.75
    "H126" = 3084.3076
    "I126" = 3116.3333
    "K126" = 9348.999899999999
    "M126" = -6878.999899999999
}
foreach ($cellRef in $updates.Keys) {
    $ws.Range($cellRef).Value = $updates[$cellRef]
}

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
$updates = @{
    "H7" = 6319.136
    "I7" = 6200.75
    "J7" = 6634.8335
    "K7" = 6200.75
    "L7" = 6634.8335
    "M7" = -6088.75
    "N7" = -6858.8335
    "H40" = 5428.5557
    "I40" = 5569.4165
    "K40" = 5569.4165
    "M40" = -5433.4165
    "H122" = 4911.1665
    "J122" = 6267.737
    "L122" = 18803.211
    "N122" = -23703.211
    "H126" = 6319.136
    "I126" = 6200.75
    "J126" = 6634.8335
    "K126" = 18602.25
    "L126" = 19904.5005
    "M126" = -16132.25
    "N126" = -24844.5005
    "H135" = 99000
    "J135" = 99000
    "L135" = 99000
    "N135" = -109140
}
foreach ($cellRef in $updates.Keys) {
    $ws.Range($cellRef).Value = $updates[$cellRef]
}

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
$updates = @{
    "H62" = 3896.6
    "I62" = 3828
    "K62" = 3828
    "M62" = -3204
    "H65" = 3896.6
    "I65" = 3828
    "K65" = 19140
    "M65" = -16020
    "H81" = 4095.25
    "I81" = 3244.1667
    "K81" = 6488.3334
    "M81" = -5427.3334
    "H84" = 4095.25
    "I84" = 3244.1667
    "K84" = 32441.667
    "M84" = -27137.667
    "H126" = 1405.7059
    "I126" = 1265.4546
    "J126" = 1662.8334
    "K126" = 3796.3638
    "L126" = 4988.5002
    "M126" = -1326.3638
    "N126" = -9928.5002
    "H132" = 2507.3333
    "I132" = 2210.5312
    "K132" = 6631.5936
    "M132" = -4101.5936
    "H136" = 1821.0769
    "I136" = 997.6774
    "J136" = 5011.75
    "K136" = 2993.0322
    "L136" = 15035.25
    "M136" = -443.0322000000001
    "N136" = -20135.25
}
foreach ($cellRef in $updates.Keys) {
    $ws.Range($cellRef).Value = $updates[$cellRef]
}
